# Update "want to go" (想去人数) counts in column F across the
# 展览 (Exhibition), 演出 (Performance) and 全部类型 (All Types) sheets.
# Values below were taken from the authoritative diff of the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$exhibitionUpdates = @{
    "F4"  = 248
    "F5"  = 1789
    "F7"  = 325
    "F8"  = 512
    "F9"  = 4677
    "F16" = 1888
    "F17" = 2986
    "F21" = 171
    "F22" = 38
    "F23" = 655
    "F24" = 933
    "F25" = 305
    "F27" = 2663
    "F28" = 1019
    "F29" = 2459
    "F30" = 249
    "F32" = 613
    "F33" = 92
    "F34" = 888
    "F35" = 428
    "F36" = 1128
    "F37" = 926
    "F38" = 1185
    "F39" = 19
    "F40" = 865
    "F42" = 363
    "F43" = 286
    "F44" = 3490
}

$performanceUpdates = @{
    "F3" = 15
}

$allTypesUpdates = @{
    "F4"  = 248
    "F6"  = 1789
    "F8"  = 325
    "F9"  = 512
    "F10" = 4677
    "F12" = 15
    "F15" = 2986
    "F20" = 171
    "F24" = 38
    "F26" = 933
    "F27" = 305
    "F28" = 2663
    "F31" = 1019
    "F32" = 2459
    "F34" = 613
    "F36" = 92
    "F37" = 888
    "F38" = 1128
    "F39" = 926
    "F41" = 1185
    "F42" = 865
    "F44" = 363
    "F48" = 286
    "F49" = 3490
}

$sheetUpdateMap = @{
    "展览"   = $exhibitionUpdates
    "演出"   = $performanceUpdates
    "全部类型" = $allTypesUpdates
}

foreach ($sheetName in $sheetUpdateMap.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $updates = $sheetUpdateMap[$sheetName]
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
